# Generate Report for handoff
# - The handoff markdown report file name changed (new GUID), and the
#   handoff status flips from "Ready for handoff" to "Handoff transform
#   failed" on every sheet. Because the transform failed, the per-locale
#   sheets lose their "Latest Handoff File" value/hyperlink and the
#   associated datetime / handback / reason columns reset to defaults.

$wb = $excel.ActiveWorkbook

$newFileName = "4430e4f1-aebe-483e-be79-9220610db35f.md"
$newStatus   = "Handoff transform failed"
$epoch       = "0001-01-01 00:00:00"
$ignored     = "Ignored"

# ---------------------------------------------------------------------
# Overview sheet: update the file-name hyperlink + its display text, and
# the status reported for each locale column.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Row -eq 2) {
        $hl.TextToDisplay = $newFileName
    }
}

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de): same filename/status update, plus the
# handoff-transform-failure resets described above.
# ---------------------------------------------------------------------
$localeSheets = @("zh-cn", "de-de")

foreach ($sheetName in $localeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A2").Value = $newFileName
    $ws.Range("B2").Value = $newStatus

    # Remove the now-stale "Latest Handoff File" hyperlink (C2) + its value.
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 3) {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime / Latest Handback DateTime reset to epoch,
    # Handoff Reason resets to "Ignored" for both data rows.
    $ws.Range("D2").Value = $epoch
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = $ignored

    $ws.Range("D3").Value = $epoch
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = $ignored

    # Keep the hyperlink display text for the report-file cell in sync.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 1) {
            $hl.TextToDisplay = $newFileName
        }
    }
}
